# Update the NATMI LR-pair sheet (Fgf16-Fgfr3) with the newly-computed TPM
# values. The "ECs" sending-cluster rows are dropped entirely (only FAPs and
# MuSCs remain as senders), and every numeric column for the surviving rows
# is refreshed with the recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 8-10 (sending cluster "ECs") - delete bottom-up so the
# remaining row numbers don't shift out from under us mid-loop.
$ws.Rows(10).Delete()
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()

# Row 2: FAPs -> Fgf16/Fgfr3 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.769244333333333
$ws.Range("H2").Value = 5.307733
$ws.Range("I2").Value = 0.8057014288865171
$ws.Range("J2").Value = 0.8057014288865172
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 8.774406859680443
$ws.Range("R2").Value = 78.969661737124
$ws.Range("S2").Value = 0.6664493367868348
$ws.Range("T2").Value = 0.6664493367868349

# Row 3: FAPs -> Fgf16/Fgfr3 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.769244333333333
$ws.Range("H3").Value = 5.307733
$ws.Range("I3").Value = 0.8057014288865171
$ws.Range("J3").Value = 0.8057014288865172
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5648773333333333
$ws.Range("N3").Value = 1.694632
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 0.9994060210284443
$ws.Range("R3").Value = 8.994654189256
$ws.Range("S3").Value = 0.07590866146813635
$ws.Range("T3").Value = 0.07590866146813637

# Row 4: FAPs -> Fgf16/Fgfr3 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.769244333333333
$ws.Range("H4").Value = 5.307733
$ws.Range("I4").Value = 0.8057014288865171
$ws.Range("J4").Value = 0.8057014288865172
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4713726666666667
$ws.Range("N4").Value = 1.414118
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 0.8339734193882221
$ws.Range("R4").Value = 7.505760774494
$ws.Range("S4").Value = 0.06334343063154599
$ws.Range("T4").Value = 0.06334343063154602

# Row 5: MuSCs -> Fgf16/Fgfr3 -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4266613333333333
$ws.Range("H5").Value = 1.279984
$ws.Range("I5").Value = 0.1942985711134829
$ws.Range("J5").Value = 0.1942985711134829
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 2.115988198705778
$ws.Range("R5").Value = 19.043893788352
$ws.Range("S5").Value = 0.1607172945394503
$ws.Range("T5").Value = 0.1607172945394503

# Row 6: MuSCs -> Fgf16/Fgfr3 -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4266613333333333
$ws.Range("H6").Value = 1.279984
$ws.Range("I6").Value = 0.1942985711134829
$ws.Range("J6").Value = 0.1942985711134829
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5648773333333333
$ws.Range("N6").Value = 1.694632
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("Q6").Value = 0.2410113162097778
$ws.Range("R6").Value = 2.169101845888
$ws.Range("S6").Value = 0.01830571962467423
$ws.Range("T6").Value = 0.01830571962467423

# Row 7: MuSCs -> Fgf16/Fgfr3 -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4266613333333333
$ws.Range("H7").Value = 1.279984
$ws.Range("I7").Value = 0.1942985711134829
$ws.Range("J7").Value = 0.1942985711134829
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4713726666666667
$ws.Range("N7").Value = 1.414118
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("Q7").Value = 0.2011164904568889
$ws.Range("R7").Value = 1.810048414112
$ws.Range("S7").Value = 0.01527555694935837
$ws.Range("T7").Value = 0.01527555694935838
